# tambah qtypack di product.
# Insert a new column before column L ("Panjang") on Sheet1 and label the
# new header cell "QtyPack" (row 2 holds the column headers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at L, shifting Panjang/Lebar/Tinggi/... etc. one
# column to the right (K stays "Lead Time", new L becomes "QtyPack").
$ws.Columns.Item(12).Insert()

# Set the new header text.
$ws.Cells.Item(2, 12).Value = "QtyPack"

# Leave the selection where the new column now sits, matching the saved
# workbook's cursor position.
$ws.Range("L2").Select()
